# Update TPM-derived NATMI metrics for the Nid1-Col13a1 ligand-receptor pair sheet
# (columns G:T, rows 2-9) to reflect the refreshed TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 100.8744713333333
$ws.Range("H2").Value = 302.623414
$ws.Range("I2").Value = 0.1452075237922473
$ws.Range("J2").Value = 0.1452075237922473
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5532856666666667
$ws.Range("N2").Value = 1.659857
$ws.Range("O2").Value = 0.4010144607159208
$ws.Range("P2").Value = 0.4010144607159208
$ws.Range("Q2").Value = 55.8123991213109
$ws.Range("R2").Value = 502.3115920917981
$ws.Range("S2").Value = 0.0582303168454423
$ws.Range("T2").Value = 0.0582303168454423
$ws.Range("G3").Value = 100.8744713333333
$ws.Range("H3").Value = 302.623414
$ws.Range("I3").Value = 0.1452075237922473
$ws.Range("J3").Value = 0.1452075237922473
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8264293333333333
$ws.Range("N3").Value = 2.479288
$ws.Range("O3").Value = 0.5989855392840792
$ws.Range("P3").Value = 0.5989855392840792
$ws.Range("Q3").Value = 83.36562209435913
$ws.Range("R3").Value = 750.2905988492321
$ws.Range("S3").Value = 0.08697720694680502
$ws.Range("T3").Value = 0.08697720694680502
$ws.Range("I4").Value = 0.7769829249672668
$ws.Range("J4").Value = 0.776982924967267
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5532856666666667
$ws.Range("N4").Value = 1.659857
$ws.Range("O4").Value = 0.4010144607159208
$ws.Range("P4").Value = 0.4010144607159208
$ws.Range("Q4").Value = 298.6434861375409
$ws.Range("R4").Value = 2687.791375237868
$ws.Range("S4").Value = 0.3115813886412273
$ws.Range("T4").Value = 0.3115813886412274
$ws.Range("I5").Value = 0.7769829249672668
$ws.Range("J5").Value = 0.776982924967267
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8264293333333333
$ws.Range("N5").Value = 2.479288
$ws.Range("O5").Value = 0.5989855392840792
$ws.Range("P5").Value = 0.5989855392840792
$ws.Range("Q5").Value = 446.0765062646791
$ws.Range("R5").Value = 4014.688556382112
$ws.Range("S5").Value = 0.4654015363260395
$ws.Range("T5").Value = 0.4654015363260396
$ws.Range("G6").Value = 53.798087
$ws.Range("H6").Value = 161.394261
$ws.Range("I6").Value = 0.07744166482137986
$ws.Range("J6").Value = 0.07744166482137986
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.5532856666666667
$ws.Range("N6").Value = 1.659857
$ws.Range("O6").Value = 0.4010144607159208
$ws.Range("P6").Value = 0.4010144607159208
$ws.Range("Q6").Value = 29.76571043118633
$ws.Range("R6").Value = 267.891393880677
$ws.Range("S6").Value = 0.03105522745528874
$ws.Range("T6").Value = 0.03105522745528874
$ws.Range("G7").Value = 53.798087
$ws.Range("H7").Value = 161.394261
$ws.Range("I7").Value = 0.07744166482137986
$ws.Range("J7").Value = 0.07744166482137986
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.8264293333333333
$ws.Range("N7").Value = 2.479288
$ws.Range("O7").Value = 0.5989855392840792
$ws.Range("P7").Value = 0.5989855392840792
$ws.Range("Q7").Value = 44.46031717401867
$ws.Range("R7").Value = 400.142854566168
$ws.Range("S7").Value = 0.04638643736609112
$ws.Range("T7").Value = 0.04638643736609112
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.2555676666666666
$ws.Range("H8").Value = 0.7667029999999999
$ws.Range("I8").Value = 0.0003678864191059829
$ws.Range("J8").Value = 0.000367886419105983
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5532856666666667
$ws.Range("N8").Value = 1.659857
$ws.Range("O8").Value = 0.4010144607159208
$ws.Range("P8").Value = 0.4010144607159208
$ws.Range("Q8").Value = 0.1414019268301111
$ws.Range("R8").Value = 1.272617341471
$ws.Range("S8").Value = 0.0001475277739624969
$ws.Range("T8").Value = 0.000147527773962497
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.2555676666666666
$ws.Range("H9").Value = 0.7667029999999999
$ws.Range("I9").Value = 0.0003678864191059829
$ws.Range("J9").Value = 0.000367886419105983
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8264293333333333
$ws.Range("N9").Value = 2.479288
$ws.Range("O9").Value = 0.5989855392840792
$ws.Range("P9").Value = 0.5989855392840792
$ws.Range("Q9").Value = 0.2112086163848889
$ws.Range("R9").Value = 1.900877547464
$ws.Range("S9").Value = 0.0002203586451434859
$ws.Range("T9").Value = 0.000220358645143486
